# Generate Report for Archive
#
# The localization status for both in-flight files flipped from
# "Ready for handoff" to "In Translation" on every sheet that surfaces a
# Status column (Overview!E:F, zh-cn!C, de-de!C). Excel's AutoFit then
# re-measured those narrower Status columns down from their old width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: Status is duplicated into the zh-cn / de-de columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- Per-locale sheets: Status lives in column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Re-fit the Status columns now that the text is shorter ---
# ColumnWidth is expressed in characters; 12.5 is the narrowest input that
# still resolves to the updated column width used for the Status columns.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
